$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-24 12:54:27"

# Column O holds the scrape "timestamp" for every data row (2-410).
# Update each one to the new timestamp recorded for this run.
$lastRow = 410
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
